# B1--and-B2-PowerPoint.pptx edit
#
# 1) The table on slide 5 (the "Type of document / Definition / Why it is
#    important" table) switches from the deck's custom table style
#    ({797459C0-ED87-48B0-9F94-03458B46754D}) to a different (built-in)
#    table style ({8EE26F08-5784-4C7E-942B-E843F71A9935}).
#
# 2) The presentation's colour theme is switched from the "Integral / Red
#    Violet" palette over to the standard "Office" palette (the deck's
#    theme part keeps the same 12 theme-colour slots, just repainted with
#    the Office colours).

$p = $ppt.ActivePresentation

# -- 1. Re-style the table on slide 5 -----------------------------------
$slide = $p.Slides.Item(5)
$tableShape = $slide.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{8EE26F08-5784-4C7E-942B-E843F71A9935}")

# -- 2. Repaint the theme colour scheme (Integral -> Office) ------------
# Order of slots in ThemeColorScheme: dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink. RGB is the standard COM long (R | G<<8 | B<<16).
$slideForTheme = $p.Slides.Item(1)
$colors = $slideForTheme.ThemeColorScheme

$colors.Item(1).RGB  = 0         # dk1      000000
$colors.Item(2).RGB  = 16777215  # lt1      FFFFFF
$colors.Item(3).RGB  = 6968388   # dk2      44546A
$colors.Item(4).RGB  = 15132391  # lt2      E7E6E6
$colors.Item(5).RGB  = 13998939  # accent1  5B9BD5
$colors.Item(6).RGB  = 3243501   # accent2  ED7D31
$colors.Item(7).RGB  = 10855845  # accent3  A5A5A5
$colors.Item(8).RGB  = 49407     # accent4  FFC000
$colors.Item(9).RGB  = 12874308  # accent5  4472C4
$colors.Item(10).RGB = 4697456   # accent6  70AD47
$colors.Item(11).RGB = 12673797  # hlink    0563C1
$colors.Item(12).RGB = 7491477   # folHlink 954F72
